$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(2, 2).Value = "112号直流"
$ws.Cells.Item(2, 3).Value = 46039.590810185182
$ws.Cells.Item(2, 4).Value = 46041.391481481478

$ws.Cells.Item(3, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(3, 2).Value = "212号直流"
$ws.Cells.Item(3, 3).Value = 46040.154317129629
$ws.Cells.Item(3, 4).Value = 46041.391481481478

$ws.Cells.Item(4, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(4, 2).Value = "111号直流"
$ws.Cells.Item(4, 3).Value = 46040.390983796293
$ws.Cells.Item(4, 4).Value = 46041.391481481478

$ws.Cells.Item(5, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(5, 2).Value = "306号直流"
$ws.Cells.Item(5, 3).Value = 46040.433900462966
$ws.Cells.Item(5, 4).Value = 46041.391481481478

$ws.Cells.Item(6, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(6, 2).Value = "109号直流"
$ws.Cells.Item(6, 3).Value = 46040.498611111114
$ws.Cells.Item(6, 4).Value = 46041.391481481478

$ws.Cells.Item(7, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(7, 2).Value = "305号直流"
$ws.Cells.Item(7, 3).Value = 46040.509837962964
$ws.Cells.Item(7, 4).Value = 46041.391481481478

$ws.Cells.Item(8, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(8, 2).Value = "107号直流"
$ws.Cells.Item(8, 3).Value = 46040.555659722224
$ws.Cells.Item(8, 4).Value = 46041.391481481478

$ws.Cells.Item(9, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(9, 2).Value = "204号直流"
$ws.Cells.Item(9, 3).Value = 46040.616168981483
$ws.Cells.Item(9, 4).Value = 46041.391481481478

$ws.Cells.Item(10, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(10, 2).Value = "105号直流"
$ws.Cells.Item(10, 3).Value = 46040.651203703703
$ws.Cells.Item(10, 4).Value = 46041.391481481478

$ws.Cells.Item(11, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(11, 2).Value = "206号直流"
$ws.Cells.Item(11, 3).Value = 46040.719317129631
$ws.Cells.Item(11, 4).Value = 46041.391481481478

$ws.Cells.Item(12, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(12, 2).Value = "308号直流"
$ws.Cells.Item(12, 3).Value = 46040.751851851855
$ws.Cells.Item(12, 4).Value = 46041.391481481478

$ws.Cells.Item(13, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(13, 2).Value = "203号直流"
$ws.Cells.Item(13, 3).Value = 46040.87871527778
$ws.Cells.Item(13, 4).Value = 46041.391481481478

$ws.Cells.Item(14, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(14, 2).Value = "302号直流"
$ws.Cells.Item(14, 3).Value = 46040.882905092592
$ws.Cells.Item(14, 4).Value = 46041.391481481478

$ws.Cells.Item(15, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(15, 2).Value = "9176699400500303"
$ws.Cells.Item(15, 3).Value = 46038.375428240739
$ws.Cells.Item(15, 4).Value = 46041.397337962961

$ws.Cells.Item(16, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(16, 2).Value = "9176699400500304"
$ws.Cells.Item(16, 3).Value = 46038.56422453704
$ws.Cells.Item(16, 4).Value = 46041.397337962961

$ws.Cells.Item(17, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(17, 2).Value = "9176699400501304"
$ws.Cells.Item(17, 3).Value = 46039.521307870367
$ws.Cells.Item(17, 4).Value = 46041.397337962961

$ws.Cells.Item(18, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(18, 2).Value = "9176699400501302"
$ws.Cells.Item(18, 3).Value = 46039.536747685182
$ws.Cells.Item(18, 4).Value = 46041.397337962961

$ws.Cells.Item(19, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(19, 2).Value = "9176699400500502"
$ws.Cells.Item(19, 3).Value = 46039.562951388885
$ws.Cells.Item(19, 4).Value = 46041.397337962961

$ws.Cells.Item(20, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(20, 2).Value = "9176699400501205"
$ws.Cells.Item(20, 3).Value = 46039.585416666669
$ws.Cells.Item(20, 4).Value = 46041.397337962961

$ws.Cells.Item(21, 1).Value = "飞狐四方坪东区充电站"
$ws.Cells.Item(21, 2).Value = "9176699420300105"
$ws.Cells.Item(21, 3).Value = 46040.035381944443
$ws.Cells.Item(21, 4).Value = 46041.397337962961

$ws.Cells.Item(22, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(22, 2).Value = "9176699400500102"
$ws.Cells.Item(22, 3).Value = 46040.035416666666
$ws.Cells.Item(22, 4).Value = 46041.397337962961

$ws.Cells.Item(23, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(23, 2).Value = "9176699355900102"
$ws.Cells.Item(23, 3).Value = 46040.062719907408
$ws.Cells.Item(23, 4).Value = 46041.397337962961

$ws.Cells.Item(24, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(24, 2).Value = "9176699400501105"
$ws.Cells.Item(24, 3).Value = 46040.080555555556
$ws.Cells.Item(24, 4).Value = 46041.397337962961

$ws.Cells.Item(25, 1).Value = "飞狐四方坪东区充电站"
$ws.Cells.Item(25, 2).Value = "9176699442100402"
$ws.Cells.Item(25, 3).Value = 46040.084699074076
$ws.Cells.Item(25, 4).Value = 46041.397337962961

$ws.Cells.Item(26, 1).Value = "飞狐四方坪东区充电站"
$ws.Cells.Item(26, 2).Value = "9176699420300104"
$ws.Cells.Item(26, 3).Value = 46040.118495370371
$ws.Cells.Item(26, 4).Value = 46041.397337962961

$ws.Cells.Item(27, 1).Value = "飞狐四方坪南区充电站"
$ws.Cells.Item(27, 2).Value = "9176699368200406"
$ws.Cells.Item(27, 3).Value = 46040.423483796294
$ws.Cells.Item(27, 4).Value = 46041.397337962961

$ws.Cells.Item(28, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(28, 2).Value = "9176699400500105"
$ws.Cells.Item(28, 3).Value = 46040.548564814817
$ws.Cells.Item(28, 4).Value = 46041.397337962961

$ws.Cells.Item(29, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(29, 2).Value = "9176699400500302"
$ws.Cells.Item(29, 3).Value = 46040.549293981479
$ws.Cells.Item(29, 4).Value = 46041.397337962961

$ws.Cells.Item(30, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(30, 2).Value = "9176699400500604"
$ws.Cells.Item(30, 3).Value = 46040.581365740742
$ws.Cells.Item(30, 4).Value = 46041.397337962961

$ws.Cells.Item(31, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(31, 2).Value = "9176699400501101"
$ws.Cells.Item(31, 3).Value = 46040.583425925928
$ws.Cells.Item(31, 4).Value = 46041.397337962961

$ws.Cells.Item(32, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(32, 2).Value = "9176699400501004"
$ws.Cells.Item(32, 3).Value = 46040.587881944448
$ws.Cells.Item(32, 4).Value = 46041.397337962961

$ws.Cells.Item(33, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(33, 2).Value = "9176699400501104"
$ws.Cells.Item(33, 3).Value = 46040.59412037037
$ws.Cells.Item(33, 4).Value = 46041.397337962961

$ws.Cells.Item(34, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(34, 2).Value = "9176699400500501"
$ws.Cells.Item(34, 3).Value = 46040.620196759257
$ws.Cells.Item(34, 4).Value = 46041.397337962961

$ws.Cells.Item(35, 1).Value = "飞狐四方坪南区充电站"
$ws.Cells.Item(35, 2).Value = "9176699368200203"
$ws.Cells.Item(35, 3).Value = 46040.624456018515
$ws.Cells.Item(35, 4).Value = 46041.397337962961

$ws.Cells.Item(36, 1).Value = "飞狐四方坪南区充电站"
$ws.Cells.Item(36, 2).Value = "9176699368200403"
$ws.Cells.Item(36, 3).Value = 46040.672013888892
$ws.Cells.Item(36, 4).Value = 46041.397337962961

$ws.Cells.Item(37, 1).Value = "飞狐四方坪东区充电站"
$ws.Cells.Item(37, 2).Value = "9176699442100301"
$ws.Cells.Item(37, 3).Value = 46040.678240740737
$ws.Cells.Item(37, 4).Value = 46041.397337962961

$ws.Cells.Item(38, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(38, 2).Value = "9176699400500205"
$ws.Cells.Item(38, 3).Value = 46040.694768518515
$ws.Cells.Item(38, 4).Value = 46041.397337962961

$ws.Cells.Item(39, 1).Value = "飞狐四方坪东区充电站"
$ws.Cells.Item(39, 2).Value = "9176699442100702"
$ws.Cells.Item(39, 3).Value = 46040.717060185183
$ws.Cells.Item(39, 4).Value = 46041.397337962961

$ws.Cells.Item(40, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(40, 2).Value = "9176699400500305"
$ws.Cells.Item(40, 3).Value = 46040.73982638889
$ws.Cells.Item(40, 4).Value = 46041.397337962961

$ws.Cells.Item(41, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(41, 2).Value = "9176699400500703"
$ws.Cells.Item(41, 3).Value = 46040.754467592589
$ws.Cells.Item(41, 4).Value = 46041.397337962961

$ws.Cells.Item(42, 1).Value = "飞狐四方坪西区充电站"
$ws.Cells.Item(42, 2).Value = "9176699400500403"
$ws.Cells.Item(42, 3).Value = 46040.79991898148
$ws.Cells.Item(42, 4).Value = 46041.397337962961

$ws.Cells.Item(43, 1).Value = ""
$ws.Cells.Item(43, 2).Value = ""
$ws.Cells.Item(43, 3).Value = ""
$ws.Cells.Item(43, 4).Value = ""

$ws.Cells.Item(44, 1).Value = ""
$ws.Cells.Item(44, 2).Value = ""
$ws.Cells.Item(44, 3).Value = ""
$ws.Cells.Item(44, 4).Value = ""
